$d = $word.ActiveDocument

function Set-RangeText($startPos, $endPos, $text) {
    # Forces a genuine text replacement (not a no-op) by first writing a
    # throwaway placeholder, which normalizes/merges the run(s) covering the
    # range, then writing the real text into the same (now single) run.
    $tmp = $d.Range($startPos, $endPos)
    $tmp.Text = [string][char]1
    $final = $d.Range($startPos, $startPos + 1)
    $final.Text = $text
    return $final
}

# ---------------------------------------------------------------------------
# 1) Bold the title paragraph ("Chapter 2_ Bipartite Matching.")
# ---------------------------------------------------------------------------
$titlePara = $d.Paragraphs.Item(1)
$titlePara.Range.Bold = 1

# ---------------------------------------------------------------------------
# 2) Paragraph: "If there are no reachable vertexes ..." - normalize the
#    three runs into one, then split it again by inserting the _GoBack
#    bookmark right after "...then it is said".
# ---------------------------------------------------------------------------
$find = $d.Content.Find
$find.ClearFormatting()
$find.Execute("If there are no reachable vertexes of the whole*traversal algorithm.", $true, $false, $true, $false, $false, $true, 1, $false, "", 0)
if ($find.Found) {
    $p = $find.Parent
    Set-RangeText $p.Start $p.End "If there are no reachable vertexes of the whole (left and right) are present in the unmatched vertex in the right then it is said to be maximum. This can be achieved by using any algorithms such as BFS and other graph traversal algorithm."
}

$find2 = $d.Content.Find
$find2.ClearFormatting()
$find2.Execute("then it is said", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($find2.Found) {
    $bmPos = $find2.Parent.End
    $bmRange = $d.Range($bmPos, $bmPos)
    $d.Bookmarks.Add("_GoBack", $bmRange)
}

# ---------------------------------------------------------------------------
# 3) Paragraph: "Let G be a bipartite graph ..." - merge the first two runs
#    (this also drops the old _GoBack bookmark that used to sit between them,
#    since the bookmark moved above).
# ---------------------------------------------------------------------------
$find3 = $d.Content.Find
$find3.ClearFormatting()
$find3.Execute("Let G be a bipartite graph consisting of sets u and w such that*every non empty set x subset of u.  \(", $true, $false, $true, $false, $false, $true, 1, $false, "", 0)
if ($find3.Found) {
    $p3 = $find3.Parent
    Set-RangeText $p3.Start $p3.End "Let G be a bipartite graph consisting of sets u and w such that |u| <= |w|. G satisfies Hall’s condition if |N(x)| >= |x| for every non empty set x subset of u.  ("
}

# ---------------------------------------------------------------------------
# 4) Paragraph: "Theorem: A bipartite graph G consisting of sets u and w, ..."
#    - merge the three runs into a single run.
# ---------------------------------------------------------------------------
$find4 = $d.Content.Find
$find4.ClearFormatting()
$find4.Execute("Theorem: A bipartite graph G consisting of sets u and w,*if it satisfies Hall’s condition", $true, $false, $true, $false, $false, $true, 1, $false, "", 0)
if ($find4.Found) {
    $p4 = $find4.Parent
    Set-RangeText $p4.Start $p4.End "Theorem: A bipartite graph G consisting of sets u and w, |u| <= |w| has a matching of size |u| (small side) if and only if it satisfies Hall’s condition"
}

Write-Output "done"
